$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add the new data row (row 4) for perssocialpsychrev ---
$ws.Range("A4").Value = 20150809
$ws.Range("B4").Value = "perssocialpsychrev"
$ws.Range("C4").Value = "http://psr.sagepub.com/content/by/year/"
$ws.Range("D4").Value = "http://psr.sagepub.com/content/by/year/[0-9]{4}"
$ws.Range("E4").Value = "http://psr.sagepub.com/content/vol[0-9]{1,}/issue[0-9]{1,}/"
$ws.Range("F4").Value = "http://psr.sagepub.com/content/[0-9]{1,}/[0-9]{1,}/[0-9]{1,}.abstract"

# --- Recreate the row-3 hyperlinks (same targets, reinserted in reverse order) ---
$urlC3 = "http://pps.sagepub.com/content/by/year/"
$urlD3 = "http://pps.sagepub.com/content/by/year/[0-9]{4}"
$urlE3 = "http://pps.sagepub.com/content/vol[0-9]{1,}/issue[0-9]{1,}/"
$urlF3 = "http://pps.sagepub.com/content/[0-9]{1,}/[0-9]{1,}/[0-9]{1,}.abstract"

$origStyle = $ws.Range("C3").Style

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F3"), $urlF3)
$ws.Hyperlinks.Add($ws.Range("E3"), $urlE3)
$ws.Hyperlinks.Add($ws.Range("D3"), $urlD3)
$ws.Hyperlinks.Add($ws.Range("C3"), $urlC3)

# Re-apply the original "Hyperlink" cell style that Add() above would otherwise
# have cloned into a brand-new style index.
$ws.Range("C3:F3").Style = $origStyle

$ws.Range("A5").Select()
